$d = $word.ActiveDocument

# 1. Replace the ID placeholder text (this also removes the trailing run that
#    only contained a single space, since Find/Replace merges the matched
#    range -- which spans both runs -- into the replacement text).
$d.Content.Find.Execute("**ID__AFFARS_5317_topic_12__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5317_205__ID**", 2)

# 2. Update the first paragraph's indentation and add a (line-less) paragraph
#    border with 5-twip spacing on all sides, matching the border later used
#    in the document.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25
$b1 = $p1.Range.ParagraphFormat.Borders
$b1.DistanceFromTop = 5
$b1.DistanceFromLeft = 5
$b1.DistanceFromBottom = 5
$b1.DistanceFromRight = 5
